# Generate Report for handback
# The 2d156d40-...md file has been handed back (in sync with en-US), so its
# status moves from "Ready for handoff" to "Handed back: in sync with en-US"
# across the Overview, zh-cn and de-de sheets, and the corresponding
# "Latest Handback DateTime" timestamps for that handback are refreshed.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = $statusText
$ov.Range("C3").Value = $statusText

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B3").Value = $statusText
$zh.Range("G2").Value = "2016-02-15 09:01:55"
$zh.Range("G3").Value = "2016-02-15 09:01:55"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("B3").Value = $statusText
$de.Range("G2").Value = "2016-02-15 09:02:23"
$de.Range("G3").Value = "2016-02-15 09:02:23"
